$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mojibake = [string]([char]194) + [string]([char]177)
$correct = [string]([char]177)

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $usedRange.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains($mojibake)) {
            $cell.Value2 = $val.Replace($mojibake, $correct)
        }
    }
}
